# 'actualización 03 01 2021'
# Append the new rows of (clave, fecha) data for 2020-12-23 .. 2020-12-31
# to the "clave_numero_fecha" sheet, continuing the existing series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Last existing data row is 201 (clave 276 / 2020-12-22). Add rows 202-207.
$startRow = 202

$codigos = @(277, 278, 279, 280, 281, 282)
$fechas  = @(44188, 44189, 44193, 44194, 44195, 44196)

# Carry the date number-format (cellXf s="1", numFmtId 14) already used by
# column B down into the newly added cells by copying the format of the
# last formatted cell before writing the new values.
$formatSource = $ws.Range("B201")

for ($i = 0; $i -lt $codigos.Length; $i++) {
    $row = $startRow + $i
    $formatSource.Copy($ws.Cells.Item($row, 2))
    $ws.Cells.Item($row, 1).Value = $codigos[$i]
    $ws.Cells.Item($row, 2).Value = $fechas[$i]
}

$lastRow = $startRow + $codigos.Length - 1
[void]$ws.Range("B" + $lastRow).Offset(1, 0).Select()
